$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 292, shifting existing rows 292:310 down to 293:311
$ws.Rows.Item(292).Insert()

# Fill in the new row 292 with the new weekly price record
$ws.Cells.Item(292, 1).Value = 4
$ws.Cells.Item(292, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(292, 3).Value = "Los Lagos"
$ws.Cells.Item(292, 4).Value = 44746
$ws.Cells.Item(292, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(292, 5).Value = 10
$ws.Cells.Item(292, 6).Value = 100112040
$ws.Cells.Item(292, 7).Value = "Cilantro"
$ws.Cells.Item(292, 8).Value = "Sin especificar"
$ws.Cells.Item(292, 9).Value = "Primera"
$ws.Cells.Item(292, 10).Value = 80
$ws.Cells.Item(292, 11).Value = 13000
$ws.Cells.Item(292, 12).Value = 13000
$ws.Cells.Item(292, 13).Value = 13000
$ws.Cells.Item(292, 14).Value = '$/caja 36 atados'
$ws.Cells.Item(292, 15).Value = "Región Metropolitana"
$ws.Cells.Item(292, 16).Value = 361
$ws.Cells.Item(292, 17).Value = 36
$ws.Cells.Item(292, 18).Value = "Hortaliza"
